$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.496.55"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").Value = "1.608.31"
$ws.Range("E3").Value = "  +2.91%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'212.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.12%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'26.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.04%  "
$ws.Range("D9").Value = "'43.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").Value = "'0.0911"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "1.837.24"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").Value = "1.630.82"
$ws.Range("E14").Value = "  +4.38%  "
$ws.Range("D15").Value = "29.505.33"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "'0.536"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "'63.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").Value = "'241.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.41%  "
$ws.Range("D20").Value = "'7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.79%  "
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("D24").Value = "'9.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("D25").Value = "'2.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "'154.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("E27").Value = "  +5.01%  "
$ws.Range("E28").Value = "  +3.36%  "
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +2.61%  "
$ws.Range("D32").Value = "'1.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +4.21%  "
$ws.Range("D35").Value = "1.413.63"
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("E38").Value = "  +5.17%  "
$ws.Range("D39").Value = "'2.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'0.0165"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.62%  "
$ws.Range("D41").Value = "'0.538"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("E43").Value = "  +5.67%  "
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'52.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +22.32%  "
$ws.Range("D47").Value = "'65.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D49").Value = "1.748.47"
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("D50").Value = "'0.862"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("E51").Value = "  +1.93%  "
